# -----------------------------------------------------------------------
# Commit: "moved strategies, added time"
#
# 1) Insert 4 new columns (Interval, Start, End, Duration) after the
#    Stop Loss [%] column and before Equity Final [$]; shifts old D:O -> H:S.
# 2) Refresh row 2's metric values (re-run of the backtest) and fill the
#    new Interval/Start/End/Duration cells.
# 3) Duplicate row 2 four more times (rows 3-6), carrying the *original*
#    (pre-edit) metric values - i.e. the previously-computed strategy runs
#    that got "moved" down while row 2 now holds the freshly re-run one.
# 4) Normalize page margins / window metadata to Excel defaults.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 4 columns for Interval / Start / End / Duration ---------
$ws.Range("D1:G1").EntireColumn.Insert()

$ws.Range("D1").Value = "Interval"
$ws.Range("E1").Value = "Start"
$ws.Range("F1").Value = "End"
$ws.Range("G1").Value = "Duration"

# --- 2. Row 2: new Interval/Start/End/Duration values -------------------
$ws.Cells.Item(2, 4).Value = "15m"

$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = "2025/03/01"
$ws.Cells.Item(2, 5).Style = "Normal"

$ws.Cells.Item(2, 6).NumberFormat = "@"
$ws.Cells.Item(2, 6).Value = "2025/03/31"
$ws.Cells.Item(2, 6).Style = "Normal"

$ws.Cells.Item(2, 7).Value = 30.98958333333333
$ws.Cells.Item(2, 7).NumberFormat = "0"

# --- Row 2: refreshed metric values (columns H..S, i.e. 8..19) ---------
$ws.Cells.Item(2, 8).Value = 1160882.784640625
$ws.Cells.Item(2, 9).Value = 1240150.96153125
$ws.Cells.Item(2, 10).Value = 81635.02685937501
$ws.Cells.Item(2, 11).Value = 16.08827846406249
$ws.Cells.Item(2, 12).Value = -2.20997640620936
$ws.Cells.Item(2, 13).Value = 19
$ws.Cells.Item(2, 14).Value = 57.89473684210527
$ws.Cells.Item(2, 15).Value = 7.713160437464572
$ws.Cells.Item(2, 16).Value = -1.498145053841626
$ws.Cells.Item(2, 17).Value = 1.003901241779004
$ws.Cells.Item(2, 18).Value = 6.427083333333333
$ws.Cells.Item(2, 19).Value = 1.427083333333333

# --- 3. Duplicate row 2 into rows 3-6 (original pre-edit metric values) -
$ws.Range("A3:A6").EntireRow.Insert()
$ws.Range("A2:S2").Copy($ws.Range("A3:S3"))
$ws.Range("A2:S2").Copy($ws.Range("A4:S4"))
$ws.Range("A2:S2").Copy($ws.Range("A5:S5"))
$ws.Range("A2:S2").Copy($ws.Range("A6:S6"))

$origRow = @(1155771.778515625, 1234551.80728125, 81291.03298437501, 15.57717785156253, `
             -2.20997640620936, 19, 57.89473684210527, 7.472265629532194, -1.48845119003933, `
             0.9731947980055722)
$rows = 3, 4, 5, 6
foreach ($r in $rows) {
    $ws.Cells.Item($r, 1).Value  = "EmaCross"
    $ws.Cells.Item($r, 8).Value  = $origRow[0]
    $ws.Cells.Item($r, 9).Value  = $origRow[1]
    $ws.Cells.Item($r, 10).Value = $origRow[2]
    $ws.Cells.Item($r, 11).Value = $origRow[3]
    $ws.Cells.Item($r, 12).Value = $origRow[4]
    $ws.Cells.Item($r, 13).Value = $origRow[5]
    $ws.Cells.Item($r, 14).Value = $origRow[6]
    $ws.Cells.Item($r, 15).Value = $origRow[7]
    $ws.Cells.Item($r, 16).Value = $origRow[8]
    $ws.Cells.Item($r, 17).Value = $origRow[9]
}

# --- 4. Cosmetic metadata: window placement + default page margins -----
$win = $wb.Windows.Item(1)
$win.Left = 240
$win.Top = 15
$win.Width = 16095
$win.Height = 9660

$ps = $ws.PageSetup
$ps.LeftMargin = 0.7 * 72
$ps.RightMargin = 0.7 * 72
$ps.TopMargin = 0.75 * 72
$ps.BottomMargin = 0.75 * 72
$ps.HeaderMargin = 0.3 * 72
$ps.FooterMargin = 0.3 * 72

$ws.Range("A1").Select()
